# "Added offer page code for admin"
#
# The workbook's first sheet (a student/candidate tracking sheet) gets:
#   1. Renamed from "Alexander, Bradley and Gonzales" to "Acosta, Butler and Perez"
#   2. A new candidate ("Inna") written into row 2
#   3. The remaining candidate rows re-shuffled / edited so the final table is:
#        row2: Inna            (new)
#        row3: Suhana Sharma   (was row 4)
#        row4: Nihar           (was row 2, Red Flags 3 -> 2)
#        row5: Sagar Shah      (was row 5, unchanged)
#      "Samay Raina" (old row 3) and "Ellen Degenerous" (old row 6) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Acosta, Butler and Perez"

function Set-RowValues {
    param($sheet, $row, $values)
    $col = 1
    foreach ($val in $values) {
        $rng = $sheet.Cells.Item($row, $col)
        # Force text storage (matches source file, where every data cell -
        # even phone numbers / dates / marks - is stored as a shared string,
        # not a number/date) without leaving a non-default cell style behind.
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = "Normal"
        $col = $col + 1
    }
}

# New row 2: Inna
Set-RowValues $ws 2 @("Inna", "2023PCP5321", "8050106439", "niharkajla123@gmail.com", "2024-01-16", "Female", "PG", "CE", "90", "90", "9", "0", "0", "SC")

# New row 3: Suhana Sharma (previously row 4)
Set-RowValues $ws 3 @("Suhana Sharma", "2023PCP5305", "8050106439", "niharamazon5005@gmail.com", "2024-03-07", "Female", "PG", "CSE", "90", "90", "9", "0", "0", "General")

# New row 4: Nihar (previously row 2; Red Flags changes from 3 to 2)
Set-RowValues $ws 4 @("Nihar", "2023PCP5317", "8050106439", "niharkajla28@gmail.com", "1995-10-28", "Male", "PG", "CSE", "88", "85", "8.167", "1", "2", "General")

# Row 5: Sagar Shah (unchanged content, rewritten for consistency)
Set-RowValues $ws 5 @("Sagar Shah", "2023PCP5319", "8050106439", "niharkajla123@gmail.com", "2024-02-01", "Male", "PG", "VLSI", "80", "85", "8.75", "0", "0", "OBC")

# Remove the old trailing row (was "Ellen Degenerous"), shrinking the table to 5 rows (header + 4).
$ws.Rows.Item(6).Delete()
